{"js": "// Section 2 (\"Status of this Document\") of the Best Practice currently\n// describes the document as a \"Draft Best Practice\" and explains that,\n// \"As a draft, it may be updated, replaced or obsoleted by other\n// documents at any time. This document should not be cited as anything\n// other than work in progress. Readers are encouraged...\"\n//\n// The commit removes the \"draft\" language: \"Draft Best Practice\" becomes\n// \"Best Practice\", and the sentence about being a draft / work in\n// progress is trimmed down to \"It may be updated, replaced or obsoleted\n// by other documents at any time. Readers are encouraged...\".\n\nconst body = context.document.body;\n\n// 1) \"Draft Best Practice\" -> \"Best Practice\"\nconst titleHits = body.search(\"Draft Best Practice\", { matchCase: true });\ntitleHits.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < titleHits.items.length; i++) {\n  titleHits.items[i].insertText(\"Best Practice\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Trim the \"draft\" sentence down, dropping the \"work in progress\" caveat.\nconst draftSentence =\n  \"As a draft, it may be updated, replaced or obsoleted by other documents \" +\n  \"at any time. This document should not be cited as anything other than \" +\n  \"work in progress. Readers are encouraged to consult the following for \" +\n  \"a list of current issues, to which they are invited to contribute.\";\nconst replacement =\n  \"It may be updated, replaced or obsoleted by other documents at any \" +\n  \"time. Readers are encouraged to consult the following for a list of \" +\n  \"current issues, to which they are invited to contribute.\";\n\nconst sentenceHits = body.search(draftSentence, { matchCase: true });\nsentenceHits.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < sentenceHits.items.length; i++) {\n  sentenceHits.items[i].insertText(replacement, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Section 2 (\"Status of this Document\") of the Best Practice currently\n# describes the document as a \"Draft Best Practice\" and explains that,\n# \"As a draft, it may be updated, replaced or obsoleted by other\n# documents at any time. This document should not be cited as anything\n# other than work in progress. Readers are encouraged...\"\n#\n# This removes the \"draft\" language: \"Draft Best Practice\" becomes\n# \"Best Practice\", and the sentence about being a draft / work in\n# progress is trimmed down to \"It may be updated, replaced or obsoleted\n# by other documents at any time. Readers are encouraged...\".\n\n$d = $word.ActiveDocument\n\n# 1) \"Draft Best Practice\" -> \"Best Practice\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Draft Best Practice\"\n$find1.Replacement.Text = \"Best Practice\"\n$find1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2) Trim the \"draft\" sentence down, dropping the \"work in progress\" caveat.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"As a draft, it may be updated, replaced or obsoleted by other documents at any time. This document should not be cited as anything other than work in progress. Readers are encouraged to consult the following for a list of current issues, to which they are invited to contribute.\"\n$find2.Replacement.Text = \"It may be updated, replaced or obsoleted by other documents at any time. Readers are encouraged to consult the following for a list of current issues, to which they are invited to contribute.\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
